# Update the cryptos price table (Sheet1) with refreshed values.
# Cells whose new value looks like a plain number (e.g. "567.33") are
# written with a leading apostrophe so Excel keeps storing them as text
# (matching how the rest of the Price/Volume columns are stored),
# instead of silently converting them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.967.92"
$ws.Range("E2").Value = "  +2.97%  "
$ws.Range("D3").Value = "2.445.57"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'567.33"
$ws.Range("E5").Value = "  +2.43%  "
$ws.Range("D6").Value = "'166.91"
$ws.Range("E6").Value = "  +4.34%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.513"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("E9").Value = "  +8.51%  "
$ws.Range("D10").Value = "2.444.65"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("E13").Value = "  -1.64%  "
$ws.Range("E14").Value = "  +6.40%  "
$ws.Range("D15").Value = "70.006.88"
$ws.Range("E15").Value = "  +3.21%  "
$ws.Range("D16").Value = "2.896.04"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "'24.09"
$ws.Range("E17").Value = "  +4.74%  "
$ws.Range("D18").Value = "2.442.55"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").Value = "'10.83"
$ws.Range("E19").Value = "  +4.49%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'341.10"
$ws.Range("E20").Value = "  +2.13%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'7.13"
$ws.Range("E21").Value = "  +4.23%  "
$ws.Range("D22").Value = "'3.89"
$ws.Range("E22").Value = "  +2.89%  "
$ws.Range("D23").Value = "'2.00"
$ws.Range("E23").Value = "  +7.74%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "'66.27"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").Value = "'3.82"
$ws.Range("D27").Value = "2.572.50"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "'8.48"
$ws.Range("E28").Value = "  +4.71%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "0.0₃0856"
$ws.Range("E30").Value = "  +5.93%  "
$ws.Range("D31").Value = "'7.38"
$ws.Range("E31").Value = "  +3.99%  "
$ws.Range("D32").Value = "'456.98"
$ws.Range("E32").Value = "  +9.52%  "
$ws.Range("E33").Value = "  +9.62%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("D36").Value = "'159.27"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").Value = "'19.09"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("E38").Value = "  +5.55%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "'18.22"
$ws.Range("D41").Value = "'0.303"
$ws.Range("E41").Value = "  +3.01%  "
$ws.Range("E42").Value = "  +4.40%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "'38.02"
$ws.Range("E43").Value = "  +1.63%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D44").Value = "'4.40"
$ws.Range("E44").Value = "  +3.38%  "
$ws.Range("E45").Value = "  +1.59%  "
$ws.Range("E46").Value = "  +5.77%  "
$ws.Range("D47").Value = "'134.69"
$ws.Range("E47").Value = "  +3.94%  "
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("E49").Value = "  +2.85%  "
$ws.Range("D50").Value = "'0.490"
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("E51").Value = "  +1.55%  "
